# ajout section tests + dictionnaire anglais pour termes
#
# The "Intergiciels" week (row 4) now also covers testing the API, so its
# "Sujet" (C4) and "Exercice" (D4) cells gain a second linked topic / are
# renamed accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "[Intergiciels](intergiciels_express.md)<br/>[Tester API](tester_api.md)"
$ws.Range("D4").Value = "[Exercice 3 - Intergiciels Express et Tests API](exercice3_intergiciels_express.md)"

# Match the author's final cursor position / scroll in the saved workbook.
$ws.Range("D5").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
